$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at position 158, pushing the existing rows 158-163
# (which become 160-165) down to make room for two new weekly records.
$ws.Rows("158:159").Insert()

# New row 158: Región Metropolitana, Primera quality record.
$ws.Cells.Item(158, 1).Value = 6
$ws.Cells.Item(158, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(158, 3).Value = "Metropolitana"
$ws.Cells.Item(158, 4).Value = 44585
$ws.Cells.Item(158, 5).Value = 13
$ws.Cells.Item(158, 6).Value = 100112001
$ws.Cells.Item(158, 7).Value = "Berenjena"
$ws.Cells.Item(158, 8).Value = "Sin especificar"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 130
$ws.Cells.Item(158, 11).Value = 9000
$ws.Cells.Item(158, 12).Value = 9000
$ws.Cells.Item(158, 13).Value = 9000
$ws.Cells.Item(158, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(158, 15).Value = "Región Metropolitana"
$ws.Cells.Item(158, 16).Value = 180
$ws.Cells.Item(158, 17).Value = 50
$ws.Cells.Item(158, 18).Value = "Hortaliza"

# New row 159: Región Metropolitana, Segunda quality record.
$ws.Cells.Item(159, 1).Value = 6
$ws.Cells.Item(159, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(159, 3).Value = "Metropolitana"
$ws.Cells.Item(159, 4).Value = 44585
$ws.Cells.Item(159, 5).Value = 13
$ws.Cells.Item(159, 6).Value = 100112001
$ws.Cells.Item(159, 7).Value = "Berenjena"
$ws.Cells.Item(159, 8).Value = "Sin especificar"
$ws.Cells.Item(159, 9).Value = "Segunda"
$ws.Cells.Item(159, 10).Value = 70
$ws.Cells.Item(159, 11).Value = 8000
$ws.Cells.Item(159, 12).Value = 8000
$ws.Cells.Item(159, 13).Value = 8000
$ws.Cells.Item(159, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(159, 15).Value = "Región Metropolitana"
$ws.Cells.Item(159, 16).Value = 160
$ws.Cells.Item(159, 17).Value = 50
$ws.Cells.Item(159, 18).Value = "Hortaliza"

# Match the date-format styling already used in column D for the new rows.
$ws.Range("D158:D159").NumberFormat = $ws.Range("D160").NumberFormat
